$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 0.1688094139099121
    3  = 0.07811379432678223
    4  = 0.1165955066680908
    5  = 0.1820666790008545
    6  = 0.301706075668335
    7  = 0.4560689926147461
    8  = 0.5821278095245361
    9  = 0.8846116065979004
    10 = 1.153199911117554
    11 = 1.496581077575684
    12 = 1.982256174087524
    13 = 2.42500376701355
    14 = 3.049225330352783
    15 = 3.665996313095093
    16 = 4.379142761230469
    17 = 5.622652530670166
    18 = 6.23734450340271
    19 = 7.259650945663452
    20 = 8.340411424636841
    21 = 9.575966835021973
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}
